# Update recomputed TPM-based NATMI ligand-receptor scores (Spon1-Lrp8)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 0.8999933333333333
$ws.Cells.Item(2, 8).Value = 2.69998
$ws.Cells.Item(2, 9).Value = 0.01781586806325543
$ws.Cells.Item(2, 10).Value = 0.01781586806325543
$ws.Cells.Item(2, 13).Value = 0.4394373333333333
$ws.Cells.Item(2, 14).Value = 1.318312
$ws.Cells.Item(2, 15).Value = 0.2944741752765458
$ws.Cells.Item(2, 16).Value = 0.2944741752765458
$ws.Cells.Item(2, 17).Value = 0.3954906704177777
$ws.Cells.Item(2, 18).Value = 3.55941603376
$ws.Cells.Item(2, 19).Value = 0.005246313054762894
$ws.Cells.Item(2, 20).Value = 0.005246313054762895

# Row 3
$ws.Cells.Item(3, 7).Value = 0.8999933333333333
$ws.Cells.Item(3, 8).Value = 2.69998
$ws.Cells.Item(3, 9).Value = 0.01781586806325543
$ws.Cells.Item(3, 10).Value = 0.01781586806325543
$ws.Cells.Item(3, 15).Value = 0.4358046333636673
$ws.Cells.Item(3, 16).Value = 0.4358046333636673
$ws.Cells.Item(3, 17).Value = 0.585303164388889
$ws.Cells.Item(3, 18).Value = 5.267728479500001
$ws.Cells.Item(3, 19).Value = 0.007764237849362502
$ws.Cells.Item(3, 20).Value = 0.007764237849362502

# Row 4
$ws.Cells.Item(4, 7).Value = 0.8999933333333333
$ws.Cells.Item(4, 8).Value = 2.69998
$ws.Cells.Item(4, 9).Value = 0.01781586806325543
$ws.Cells.Item(4, 10).Value = 0.01781586806325543
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.1607546666666667
$ws.Cells.Item(4, 14).Value = 0.482264
$ws.Cells.Item(4, 15).Value = 0.1077243426939663
$ws.Cells.Item(4, 16).Value = 0.1077243426939663
$ws.Cells.Item(4, 17).Value = 0.1446781283022222
$ws.Cells.Item(4, 18).Value = 1.30210315472
$ws.Cells.Item(4, 19).Value = 0.001919202676636618
$ws.Cells.Item(4, 20).Value = 0.001919202676636618

# Row 5
$ws.Cells.Item(5, 7).Value = 0.8999933333333333
$ws.Cells.Item(5, 8).Value = 2.69998
$ws.Cells.Item(5, 9).Value = 0.01781586806325543
$ws.Cells.Item(5, 10).Value = 0.01781586806325543
$ws.Cells.Item(5, 13).Value = 0.2417443333333333
$ws.Cells.Item(5, 14).Value = 0.725233
$ws.Cells.Item(5, 15).Value = 0.1619968486658205
$ws.Cells.Item(5, 16).Value = 0.1619968486658205
$ws.Cells.Item(5, 17).Value = 0.2175682883711111
$ws.Cells.Item(5, 18).Value = 1.95811459534
$ws.Cells.Item(5, 19).Value = 0.002886114482493415
$ws.Cells.Item(5, 20).Value = 0.002886114482493415

# Row 6
$ws.Cells.Item(6, 9).Value = 0.7949938412397365
$ws.Cells.Item(6, 10).Value = 0.7949938412397366
$ws.Cells.Item(6, 13).Value = 0.4394373333333333
$ws.Cells.Item(6, 14).Value = 1.318312
$ws.Cells.Item(6, 15).Value = 0.2944741752765458
$ws.Cells.Item(6, 16).Value = 0.2944741752765458
$ws.Cells.Item(6, 17).Value = 17.64789939696356
$ws.Cells.Item(6, 18).Value = 158.831094572672
$ws.Cells.Item(6, 19).Value = 0.2341051557490046
$ws.Cells.Item(6, 20).Value = 0.2341051557490046

# Row 7
$ws.Cells.Item(7, 9).Value = 0.7949938412397365
$ws.Cells.Item(7, 10).Value = 0.7949938412397366
$ws.Cells.Item(7, 15).Value = 0.4358046333636673
$ws.Cells.Item(7, 16).Value = 0.4358046333636673
$ws.Cells.Item(7, 19).Value = 0.3464619995078569
$ws.Cells.Item(7, 20).Value = 0.3464619995078569

# Row 8
$ws.Cells.Item(8, 9).Value = 0.7949938412397365
$ws.Cells.Item(8, 10).Value = 0.7949938412397366
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 0.1607546666666667
$ws.Cells.Item(8, 14).Value = 0.482264
$ws.Cells.Item(8, 15).Value = 0.1077243426939663
$ws.Cells.Item(8, 16).Value = 0.1077243426939663
$ws.Cells.Item(8, 17).Value = 6.455942565020446
$ws.Cells.Item(8, 18).Value = 58.10348308518401
$ws.Cells.Item(8, 19).Value = 0.08564018899330203
$ws.Cells.Item(8, 20).Value = 0.08564018899330203

# Row 9
$ws.Cells.Item(9, 9).Value = 0.7949938412397365
$ws.Cells.Item(9, 10).Value = 0.7949938412397366
$ws.Cells.Item(9, 13).Value = 0.2417443333333333
$ws.Cells.Item(9, 14).Value = 0.725233
$ws.Cells.Item(9, 15).Value = 0.1619968486658205
$ws.Cells.Item(9, 16).Value = 0.1619968486658205
$ws.Cells.Item(9, 17).Value = 9.708505288094225
$ws.Cells.Item(9, 18).Value = 87.37654759284801
$ws.Cells.Item(9, 19).Value = 0.128786496989573
$ws.Cells.Item(9, 20).Value = 0.128786496989573

# Row 10
$ws.Cells.Item(10, 7).Value = 5.293300666666666
$ws.Cells.Item(10, 8).Value = 15.879902
$ws.Cells.Item(10, 9).Value = 0.104783827617029
$ws.Cells.Item(10, 10).Value = 0.104783827617029
$ws.Cells.Item(10, 13).Value = 0.4394373333333333
$ws.Cells.Item(10, 14).Value = 1.318312
$ws.Cells.Item(10, 15).Value = 0.2944741752765458
$ws.Cells.Item(10, 16).Value = 0.2944741752765458
$ws.Cells.Item(10, 17).Value = 2.326073929491555
$ws.Cells.Item(10, 18).Value = 20.934665365424
$ws.Cells.Item(10, 19).Value = 0.03085613121984437
$ws.Cells.Item(10, 20).Value = 0.03085613121984437

# Row 11
$ws.Cells.Item(11, 7).Value = 5.293300666666666
$ws.Cells.Item(11, 8).Value = 15.879902
$ws.Cells.Item(11, 9).Value = 0.104783827617029
$ws.Cells.Item(11, 10).Value = 0.104783827617029
$ws.Cells.Item(11, 15).Value = 0.4358046333636673
$ws.Cells.Item(11, 16).Value = 0.4358046333636673
$ws.Cells.Item(11, 17).Value = 3.442453977727778
$ws.Cells.Item(11, 18).Value = 30.98208579955
$ws.Cells.Item(11, 19).Value = 0.04566527757708105
$ws.Cells.Item(11, 20).Value = 0.04566527757708105

# Row 12
$ws.Cells.Item(12, 7).Value = 5.293300666666666
$ws.Cells.Item(12, 8).Value = 15.879902
$ws.Cells.Item(12, 9).Value = 0.104783827617029
$ws.Cells.Item(12, 10).Value = 0.104783827617029
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 0.1607546666666667
$ws.Cells.Item(12, 14).Value = 0.482264
$ws.Cells.Item(12, 15).Value = 0.1077243426939663
$ws.Cells.Item(12, 16).Value = 0.1077243426939663
$ws.Cells.Item(12, 17).Value = 0.8509227842364445
$ws.Cells.Item(12, 18).Value = 7.658305058128
$ws.Cells.Item(12, 19).Value = 0.01128776895500233
$ws.Cells.Item(12, 20).Value = 0.01128776895500232

# Row 13
$ws.Cells.Item(13, 7).Value = 5.293300666666666
$ws.Cells.Item(13, 8).Value = 15.879902
$ws.Cells.Item(13, 9).Value = 0.104783827617029
$ws.Cells.Item(13, 10).Value = 0.104783827617029
$ws.Cells.Item(13, 13).Value = 0.2417443333333333
$ws.Cells.Item(13, 14).Value = 0.725233
$ws.Cells.Item(13, 15).Value = 0.1619968486658205
$ws.Cells.Item(13, 16).Value = 0.1619968486658205
$ws.Cells.Item(13, 17).Value = 1.279625440796222
$ws.Cells.Item(13, 18).Value = 11.516628967166
$ws.Cells.Item(13, 19).Value = 0.01697464986510127
$ws.Cells.Item(13, 20).Value = 0.01697464986510127

# Row 14
$ws.Cells.Item(14, 7).Value = 4.162877000000001
$ws.Cells.Item(14, 8).Value = 12.488631
$ws.Cells.Item(14, 9).Value = 0.08240646307997901
$ws.Cells.Item(14, 10).Value = 0.08240646307997902
$ws.Cells.Item(14, 13).Value = 0.4394373333333333
$ws.Cells.Item(14, 14).Value = 1.318312
$ws.Cells.Item(14, 15).Value = 0.2944741752765458
$ws.Cells.Item(14, 16).Value = 0.2944741752765458
$ws.Cells.Item(14, 17).Value = 1.829323567874667
$ws.Cells.Item(14, 18).Value = 16.463912110872
$ws.Cells.Item(14, 19).Value = 0.02426657525293394
$ws.Cells.Item(14, 20).Value = 0.02426657525293395

# Row 15
$ws.Cells.Item(15, 7).Value = 4.162877000000001
$ws.Cells.Item(15, 8).Value = 12.488631
$ws.Cells.Item(15, 9).Value = 0.08240646307997901
$ws.Cells.Item(15, 10).Value = 0.08240646307997902
$ws.Cells.Item(15, 15).Value = 0.4358046333636673
$ws.Cells.Item(15, 16).Value = 0.4358046333636673
$ws.Cells.Item(15, 17).Value = 2.707292366308334
$ws.Cells.Item(15, 18).Value = 24.36563129677501
$ws.Cells.Item(15, 19).Value = 0.03591311842936684
$ws.Cells.Item(15, 20).Value = 0.03591311842936684

# Row 16
$ws.Cells.Item(16, 7).Value = 4.162877000000001
$ws.Cells.Item(16, 8).Value = 12.488631
$ws.Cells.Item(16, 9).Value = 0.08240646307997901
$ws.Cells.Item(16, 10).Value = 0.08240646307997902
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 0.1607546666666667
$ws.Cells.Item(16, 14).Value = 0.482264
$ws.Cells.Item(16, 15).Value = 0.1077243426939663
$ws.Cells.Item(16, 16).Value = 0.1077243426939663
$ws.Cells.Item(16, 17).Value = 0.6692019045093336
$ws.Cells.Item(16, 18).Value = 6.022817140584001
$ws.Cells.Item(16, 19).Value = 0.008877182069025344
$ws.Cells.Item(16, 20).Value = 0.008877182069025344

# Row 17
$ws.Cells.Item(17, 7).Value = 4.162877000000001
$ws.Cells.Item(17, 8).Value = 12.488631
$ws.Cells.Item(17, 9).Value = 0.08240646307997901
$ws.Cells.Item(17, 10).Value = 0.08240646307997902
$ws.Cells.Item(17, 13).Value = 0.2417443333333333
$ws.Cells.Item(17, 14).Value = 0.725233
$ws.Cells.Item(17, 15).Value = 0.1619968486658205
$ws.Cells.Item(17, 16).Value = 0.1619968486658205
$ws.Cells.Item(17, 17).Value = 1.006351925113667
$ws.Cells.Item(17, 18).Value = 9.057167326023002
$ws.Cells.Item(17, 19).Value = 0.01334958732865289
$ws.Cells.Item(17, 20).Value = 0.01334958732865289

